$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('94154', 'Bleach', '2', '$15.34', '$30.68'),
    @('11387', 'Veggie Burger', '2', '$52.72', '$105.44'),
    @('11072', 'Vegan Chicken Tenders', '1', '$87.80', '$87.80'),
    @('32017', 'Sausage - Chicken Patty', '4', '$51.14', '$204.56'),
    @('3357597', 'FRZ Fruit - Strawberry (Whole)', '3', '$43.74', '$131.22'),
    @('45574', 'FRZ Fruit - Raspberry (Whole)', '1', '$35.17', '$35.17'),
    @('1345057', 'FRZ Fruit - Blueberry', '1', '$76.72', '$76.72'),
    @('45528', 'FRZ Fruit - Pineapple (Chunk)', '1', '$49.78', '$49.78'),
    @('46209', 'FRZ Banana - Slices', '1', '$21.78', '$21.78'),
    @('1138748', 'Nirvana Water - 1.5 Liter', '2', '$9.59', '$19.18'),
    @('1667864', 'Employee Water', '3', '$4.73', '$14.19'),
    @('688660', 'Hals - Black Cherry', '1', '$18.35', '$18.35'),
    @('6435380', 'Hals - Lime', '1', '$18.35', '$18.35')
)

$startRow = 9
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
    }
}
